$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 57: Quit Your Jib-jab / Gold Needle
$ws.Range("H57").Value = 33799.6
$ws.Range("J57").Value = 33799.6
$ws.Range("L57").Value = 101398.8
$ws.Range("N57").Value = -102396.8

# Row 69: Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws.Range("H69").Value = 4104.1304
$ws.Range("I69").Value = 1221.6666
$ws.Range("K69").Value = 3664.9998
$ws.Range("M69").Value = -2790.9998

# Row 72: Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws.Range("H72").Value = 4104.1304
$ws.Range("I72").Value = 1221.6666
$ws.Range("K72").Value = 10994.9994
$ws.Range("M72").Value = -6626.999400000001

# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 73810870
$ws.Range("J100").Value = 166668820
$ws.Range("L100").Value = 166668820
$ws.Range("N100").Value = -166669902

# Row 111: An Eye for Healing / Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 2789
$ws.Range("I111").Value = 2789
$ws.Range("K111").Value = 8367
$ws.Range("M111").Value = -5300

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 73194.36
$ws.Range("J112").Value = 73194.36
$ws.Range("L112").Value = 219583.08
$ws.Range("N112").Value = -221799.08

# Row 133: Big Brush, Big Dreams / Ginseng Angle Brush
$ws.Range("H133").Value = 123119.555
$ws.Range("J133").Value = 123119.555
$ws.Range("L133").Value = 123119.555
$ws.Range("N133").Value = -133239.555

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 5557630.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 5557630.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 16672891.5
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -16677991.5

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 5033.4546
$ws.Range("I138").Value = 11407.467
$ws.Range("J138").Value = 3158.745
$ws.Range("K138").Value = 34222.401
$ws.Range("L138").Value = 9476.235000000001
$ws.Range("M138").Value = -29082.401
$ws.Range("N138").Value = -19756.235

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 95304.91
$ws.Range("I45").Value = 170987.5
$ws.Range("J45").Value = 4485.8
$ws.Range("K45").Value = 170987.5
$ws.Range("L45").Value = 4485.8
$ws.Range("M45").Value = -170610.5
$ws.Range("N45").Value = -5239.8

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2358964
$ws.Range("I61").Value = 5928.5
$ws.Range("J61").Value = 11771106
$ws.Range("K61").Value = 5928.5
$ws.Range("L61").Value = 11771106
$ws.Range("M61").Value = -5716.5
$ws.Range("N61").Value = -11771530

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 930743.7
$ws.Range("I74").Value = 1863.7
$ws.Range("J74").Value = 1594229.4
$ws.Range("K74").Value = 1863.7
$ws.Range("L74").Value = 1594229.4
$ws.Range("M74").Value = -989.7
$ws.Range("N74").Value = -1595977.4

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 930743.7
$ws.Range("I77").Value = 1863.7
$ws.Range("J77").Value = 1594229.4
$ws.Range("K77").Value = 9318.5
$ws.Range("L77").Value = 7971147
$ws.Range("M77").Value = -4950.5
$ws.Range("N77").Value = -7979883

# Row 97: Ore for Me / High Steel Ingot
$ws.Range("H97").Value = 6967.684
$ws.Range("I97").Value = 8215.200000000001
$ws.Range("K97").Value = 8215.200000000001
$ws.Range("M97").Value = -7719.200000000001

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 3313
$ws.Range("I102").Value = 3313
$ws.Range("K102").Value = 3313
$ws.Range("M102").Value = -1691

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 4014
$ws.Range("I122").Value = 3649
$ws.Range("K122").Value = 10947
$ws.Range("M122").Value = -8497

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 3070.7368
$ws.Range("I132").Value = 1837.0555
$ws.Range("K132").Value = 5511.166499999999
$ws.Range("M132").Value = -2981.166499999999

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2358964
$ws.Range("I136").Value = 5928.5
$ws.Range("J136").Value = 11771106
$ws.Range("K136").Value = 17785.5
$ws.Range("L136").Value = 35313318
$ws.Range("M136").Value = -15235.5
$ws.Range("N136").Value = -35318418

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 1407.5
$ws.Range("I20").Value = 1237.7
$ws.Range("K20").Value = 1237.7
$ws.Range("M20").Value = -990.7

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 14249.186
$ws.Range("I107").Value = 16111.429
$ws.Range("J107").Value = 7731.3335
$ws.Range("K107").Value = 16111.429
$ws.Range("L107").Value = 7731.3335
$ws.Range("M107").Value = -14191.429
$ws.Range("N107").Value = -11571.3335

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 25716554
$ws.Range("I134").Value = 1867.1724
$ws.Range("J134").Value = 150004210
$ws.Range("K134").Value = 5601.5172
$ws.Range("L134").Value = 450012630
$ws.Range("M134").Value = -3066.5172
$ws.Range("N134").Value = -450017700

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 4502.769
$ws.Range("I16").Value = 3242.7144
$ws.Range("J16").Value = 5972.8335
$ws.Range("K16").Value = 3242.7144
$ws.Range("L16").Value = 5972.8335
$ws.Range("M16").Value = -2955.7144
$ws.Range("N16").Value = -6546.8335

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3030.2295
$ws.Range("I31").Value = 1857.909
$ws.Range("J31").Value = 3288.14
$ws.Range("K31").Value = 1857.909
$ws.Range("L31").Value = 3288.14
$ws.Range("M31").Value = -1562.909
$ws.Range("N31").Value = -3878.14

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3030.2295
$ws.Range("I34").Value = 1857.909
$ws.Range("J34").Value = 3288.14
$ws.Range("K34").Value = 1857.909
$ws.Range("L34").Value = 3288.14
$ws.Range("M34").Value = -1655.909
$ws.Range("N34").Value = -3692.14

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 4502.769
$ws.Range("I113").Value = 3242.7144
$ws.Range("J113").Value = 5972.8335
$ws.Range("K113").Value = 3242.7144
$ws.Range("L113").Value = 5972.8335
$ws.Range("M113").Value = -1072.7144
$ws.Range("N113").Value = -10312.8335

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 30307142
$ws.Range("I132").Value = 4714.5
$ws.Range("J132").Value = 111113620
$ws.Range("K132").Value = 14143.5
$ws.Range("L132").Value = 333340860
$ws.Range("M132").Value = -11613.5
$ws.Range("N132").Value = -333345920

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 2712.3333
$ws.Range("I134").Value = 2599.6667
$ws.Range("K134").Value = 7799.000100000001
$ws.Range("M134").Value = -5264.000100000001

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food / Table Salt
$ws.Range("H2").Value = 526.8570999999999
$ws.Range("I2").Value = 330.22223
$ws.Range("J2").Value = 674.3333
$ws.Range("K2").Value = 1981.33338
$ws.Range("L2").Value = 4045.9998
$ws.Range("M2").Value = -1868.33338
$ws.Range("N2").Value = -4271.9998

# Row 4: In Hot Water / Boiled Egg
$ws.Range("H4").Value = 1151651.2
$ws.Range("I4").Value = 1156390.2
$ws.Range("K4").Value = 3469170.6
$ws.Range("M4").Value = -3469058.6

# Row 34: Fever Pitch / Chamomile Tea
$ws.Range("H34").Value = 240.69565
$ws.Range("I34").Value = 111.5
$ws.Range("J34").Value = 253
$ws.Range("K34").Value = 334.5
$ws.Range("L34").Value = 759
$ws.Range("M34").Value = -250.5
$ws.Range("N34").Value = -927

# Row 58: Bread in the Clouds / La Noscean Toast
$ws.Range("H58").Value = 18750
$ws.Range("I58").Value = 15000
$ws.Range("K58").Value = 45000
$ws.Range("M58").Value = -44872

# Row 112: Sweet Tooth / Caramels
$ws.Range("H112").Value = 10999.857
$ws.Range("J112").Value = 14333
$ws.Range("L112").Value = 42999
$ws.Range("N112").Value = -45215

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 1810.8125
$ws.Range("I132").Value = 1775
$ws.Range("J132").Value = 1889.6
$ws.Range("K132").Value = 15975
$ws.Range("L132").Value = 17006.4
$ws.Range("M132").Value = -13445
$ws.Range("N132").Value = -22066.4

# Row 140: Sweet, Sweet Bean Juice / Mesquite Juice
$ws.Range("H140").Value = 31252440
$ws.Range("I140").Value = 35716076
$ws.Range("K140").Value = 107148228
$ws.Range("M140").Value = -107143048

# Row 141: Ocean Explosion / Acqua Pazza
$ws.Range("H141").Value = 2533.8635
$ws.Range("J141").Value = 5155
$ws.Range("L141").Value = 15465
$ws.Range("N141").Value = -25825

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers / Copper Ingot
$ws.Range("H2").Value = 4102.8
$ws.Range("I2").Value = 73.916664
$ws.Range("J2").Value = 7821.769
$ws.Range("K2").Value = 73.916664
$ws.Range("L2").Value = 7821.769
$ws.Range("M2").Value = 39.083336
$ws.Range("N2").Value = -8047.769

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 5980.3076
$ws.Range("I122").Value = 2161.1667
$ws.Range("J122").Value = 9253.857
$ws.Range("K122").Value = 6483.500100000001
$ws.Range("L122").Value = 27761.571
$ws.Range("M122").Value = -4033.500100000001
$ws.Range("N122").Value = -32661.571

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 6234.45
$ws.Range("J46").Value = 1492.6428
$ws.Range("L46").Value = 1492.6428
$ws.Range("N46").Value = -1868.6428

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 12504477
$ws.Range("I61").Value = 16671203
$ws.Range("J61").Value = 4298.75
$ws.Range("K61").Value = 16671203
$ws.Range("L61").Value = 4298.75
$ws.Range("M61").Value = -16671001
$ws.Range("N61").Value = -4702.75

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 949.3570999999999
$ws.Range("I93").Value = 931.6842
$ws.Range("J93").Value = 986.6667
$ws.Range("K93").Value = 931.6842
$ws.Range("L93").Value = 986.6667
$ws.Range("M93").Value = 316.3158
$ws.Range("N93").Value = -3482.6667

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 12504477
$ws.Range("I113").Value = 16671203
$ws.Range("J113").Value = 4298.75
$ws.Range("K113").Value = 16671203
$ws.Range("L113").Value = 4298.75
$ws.Range("M113").Value = -16669033
$ws.Range("N113").Value = -8638.75

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 10581.833
$ws.Range("I132").Value = 4357
$ws.Range("K132").Value = 13071
$ws.Range("M132").Value = -10541

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 5650.2666
$ws.Range("I62").Value = 5423.857
$ws.Range("K62").Value = 5423.857
$ws.Range("M62").Value = -4799.857

# Row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 5650.2666
$ws.Range("I65").Value = 5423.857
$ws.Range("K65").Value = 27119.285
$ws.Range("M65").Value = -23999.285

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 79380.46000000001
$ws.Range("I132").Value = 102445.1
$ws.Range("J132").Value = 2498.3333
$ws.Range("K132").Value = 307335.3
$ws.Range("L132").Value = 7494.999899999999
$ws.Range("M132").Value = -304805.3
$ws.Range("N132").Value = -12554.9999

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 33878.84
$ws.Range("I136").Value = 51076.1
$ws.Range("K136").Value = 153228.3
$ws.Range("M136").Value = -150678.3
